# Trade #4 closed at 2026-02-17 07:52:31 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" sheets with the aggregate
# effect of a new (4th) MarketMaking trade, and appends the new trade row
# to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.93             # Current Capital
$summary.Range("B4").Value = -0.07000000000000001 # Total P&L $
$summary.Range("B5").Value = -0.35               # Total P&L %
$summary.Range("B6").Value = 4                   # Total Trades
$summary.Range("B8").Value = 2                   # Losing Trades
$summary.Range("B9").Value = 50                  # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.93000000000001    # Capital
$status.Range("D4").Value = 4                    # Trades
$status.Range("E4").Value = -0.07000000000000001 # P&L $
$status.Range("F4").Value = -0.07000000000000001 # P&L %
$status.Range("G4").Value = 50                   # Win Rate %

# ---------------------------------------------------------------
# Append the new trade (#4) to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 5

    $ws.Cells.Item($row, 1).Value = 4

    # Force the date-looking text to stay text instead of being
    # auto-converted into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "07:52:25"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.62
    $ws.Cells.Item($row, 7).Value = 0.61
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -1.6129
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 99.93000000000001
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}
